$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 19609928
$ws.Range("I92").Value = 27779328
$ws.Range("J92").Value = 3366.2
$ws.Range("K92").Value = 27779328
$ws.Range("L92").Value = 3366.2
$ws.Range("M92").Value = -27778080
$ws.Range("N92").Value = -5862.2
$ws.Range("H132").Value = 2941.5527
$ws.Range("I132").Value = 2799.9722
$ws.Range("J132").Value = 5490
$ws.Range("K132").Value = 8399.9166
$ws.Range("L132").Value = 16470
$ws.Range("M132").Value = -5869.9166
$ws.Range("N132").Value = -21530
$ws.Range("H138").Value = 4235.7593
$ws.Range("I138").Value = 3744.3157
$ws.Range("J138").Value = 4391.3833
$ws.Range("K138").Value = 11232.9471
$ws.Range("L138").Value = 13174.1499
$ws.Range("M138").Value = -6092.947100000001
$ws.Range("N138").Value = -23454.1499
$ws.Range("H140").Value = 77106
$ws.Range("J140").Value = 77106
$ws.Range("L140").Value = 77106
$ws.Range("N140").Value = -87466

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15644543
$ws.Range("I32").Value = 21293452
$ws.Range("J32").Value = 26970.588
$ws.Range("K32").Value = 21293452
$ws.Range("L32").Value = 26970.588
$ws.Range("M32").Value = -21293165
$ws.Range("N32").Value = -27544.588
$ws.Range("H106").Value = 62500
$ws.Range("J106").Value = 62500
$ws.Range("L106").Value = 62500
$ws.Range("N106").Value = -65024

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 30304736
$ws.Range("I20").Value = 1478.0416
$ws.Range("J20").Value = 111113420
$ws.Range("K20").Value = 1478.0416
$ws.Range("L20").Value = 111113420
$ws.Range("M20").Value = -1231.0416
$ws.Range("N20").Value = -111113914
$ws.Range("H94").Value = 91315.37
$ws.Range("I94").Value = 100356.9
$ws.Range("J94").Value = 900
$ws.Range("K94").Value = 100356.9
$ws.Range("L94").Value = 900
$ws.Range("M94").Value = -99905.89999999999
$ws.Range("N94").Value = -1802
$ws.Range("H105").Value = 20835144
$ws.Range("I105").Value = 20835144
$ws.Range("K105").Value = 20835144
$ws.Range("M105").Value = -20833397
$ws.Range("H107").Value = 84404.914
$ws.Range("I107").Value = 111929.89
$ws.Range("J107").Value = 1830
$ws.Range("K107").Value = 111929.89
$ws.Range("L107").Value = 1830
$ws.Range("M107").Value = -110009.89
$ws.Range("N107").Value = -5670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 78624
$ws.Range("J96").Value = 78624
$ws.Range("L96").Value = 78624
$ws.Range("N96").Value = -84116
$ws.Range("H107").Value = 6945412.5
$ws.Range("I107").Value = 10417386
$ws.Range("K107").Value = 10417386
$ws.Range("M107").Value = -10415466
$ws.Range("H110").Value = 68702
$ws.Range("J110").Value = 68702
$ws.Range("L110").Value = 68702
$ws.Range("N110").Value = -76882
$ws.Range("H111").Value = 36666.5
$ws.Range("J111").Value = 36666.5
$ws.Range("L111").Value = 36666.5
$ws.Range("N111").Value = -44846.5
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 54.590908
$ws.Range("I12").Value = 36.235294
$ws.Range("J12").Value = 117
$ws.Range("K12").Value = 108.705882
$ws.Range("L12").Value = 351
$ws.Range("M12").Value = 64.294118
$ws.Range("N12").Value = -697
$ws.Range("H80").Value = 5299.625
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 5326.864
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 15980.592
$ws.Range("M80").Value = -14064
$ws.Range("N80").Value = -17852.592
$ws.Range("H83").Value = 5299.625
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 5326.864
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 47941.776
$ws.Range("M83").Value = -40320
$ws.Range("N83").Value = -57301.776
$ws.Range("H117").Value = 2827.3845
$ws.Range("I117").Value = 2439
$ws.Range("K117").Value = 7317
$ws.Range("M117").Value = -3875
$ws.Range("H121").Value = 838.1111
$ws.Range("I121").Value = 412.5
$ws.Range("J121").Value = 1178.6
$ws.Range("K121").Value = 1237.5
$ws.Range("L121").Value = 3535.8
$ws.Range("M121").Value = 72.5
$ws.Range("N121").Value = -6155.799999999999
$ws.Range("H122").Value = 2596.0862
$ws.Range("I122").Value = 424.2903
$ws.Range("J122").Value = 5089.6294
$ws.Range("K122").Value = 3818.6127
$ws.Range("L122").Value = 45806.6646
$ws.Range("M122").Value = -1368.6127
$ws.Range("N122").Value = -50706.6646
$ws.Range("H129").Value = 1596667.4
$ws.Range("J129").Value = 2166702.5
$ws.Range("L129").Value = 6500107.5
$ws.Range("N129").Value = -6510107.5
$ws.Range("H131").Value = 2990.6333
$ws.Range("J131").Value = 3505.7144
$ws.Range("L131").Value = 10517.1432
$ws.Range("N131").Value = -20597.1432
$ws.Range("H132").Value = 2318.6948
$ws.Range("I132").Value = 2525.037
$ws.Range("J132").Value = 2144.5938
$ws.Range("K132").Value = 22725.333
$ws.Range("L132").Value = 19301.3442
$ws.Range("M132").Value = -20195.333
$ws.Range("N132").Value = -24361.3442

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 394
$ws.Range("I107").Value = 345.3
$ws.Range("J107").Value = 556.3333
$ws.Range("K107").Value = 345.3
$ws.Range("L107").Value = 556.3333
$ws.Range("M107").Value = 1574.7
$ws.Range("N107").Value = -4396.3333
$ws.Range("H132").Value = 34489284
$ws.Range("I132").Value = 58831876
$ws.Range("K132").Value = 176495628
$ws.Range("M132").Value = -176493098

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 38717.152
$ws.Range("I100").Value = 73109.08
$ws.Range("J100").Value = 4325.231
$ws.Range("K100").Value = 73109.08
$ws.Range("L100").Value = 4325.231
$ws.Range("M100").Value = -72568.08
$ws.Range("N100").Value = -5407.231
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 77420
$ws.Range("J119").Value = 77420
$ws.Range("L119").Value = 77420
$ws.Range("N119").Value = -87096
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 5879.971
$ws.Range("J122").Value = 7299.9473
$ws.Range("L122").Value = 21899.8419
$ws.Range("N122").Value = -26799.8419
$ws.Range("H123").Value = 46053.625
$ws.Range("J123").Value = 46053.625
$ws.Range("L123").Value = 46053.625
$ws.Range("N123").Value = -55853.625
$ws.Range("H130").Value = 68180
$ws.Range("J130").Value = 68180
$ws.Range("L130").Value = 68180
$ws.Range("N130").Value = -78220

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 35065
$ws.Range("I75").Value = 20000
$ws.Range("J75").Value = 50130
$ws.Range("K75").Value = 20000
$ws.Range("L75").Value = 50130
$ws.Range("M75").Value = -19064
$ws.Range("N75").Value = -52002
$ws.Range("H76").Value = 94586.5
$ws.Range("J76").Value = 94586.5
$ws.Range("L76").Value = 94586.5
$ws.Range("N76").Value = -95216.5
$ws.Range("H78").Value = 35065
$ws.Range("I78").Value = 20000
$ws.Range("J78").Value = 50130
$ws.Range("K78").Value = 60000
$ws.Range("L78").Value = 150390
$ws.Range("M78").Value = -55320
$ws.Range("N78").Value = -159750
$ws.Range("H79").Value = 94586.5
$ws.Range("J79").Value = 94586.5
$ws.Range("L79").Value = 94586.5
$ws.Range("N79").Value = -96770.5
$ws.Range("H132").Value = 4168378.5
$ws.Range("I132").Value = 1558.1086
$ws.Range("J132").Value = 12154785
$ws.Range("K132").Value = 4674.325800000001
$ws.Range("L132").Value = 36464355
$ws.Range("M132").Value = -2144.325800000001
$ws.Range("N132").Value = -36469415
